$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.477.74'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '2.577.71'
$ws.Range('E3').Value = '  -2.62%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '541.57'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').Value = '143.94'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.580'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').Value = '6.76'
$ws.Range('E9').Value = '  +1.11%  '
$ws.Range('D10').Value = '0.100'
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('E11').Value = '  +3.03%  '
$ws.Range('D12').Value = '0.332'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').Value = '3.027.56'
$ws.Range('E13').Value = '  -2.96%  '
$ws.Range('D14').Value = '58.400.65'
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = '20.59'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').Value = '2.576.84'
$ws.Range('E16').Value = '  -3.27%  '
$ws.Range('D17').Value = '0.0000132'
$ws.Range('E17').Value = '  -2.73%  '
$ws.Range('D18').Value = '4.47'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('D19').Value = '334.71'
$ws.Range('E19').Value = '  -2.95%  '
$ws.Range('D20').Value = '10.03'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').Value = '6.11'
$ws.Range('E21').Value = '  -3.83%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = '66.45'
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').Value = '0.422'
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -4.92%  '
$ws.Range('D27').Value = '7.05'
$ws.Range('E27').Value = '  -3.35%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0734'
$ws.Range('E28').Value = '  -2.66%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -0.99%  '
$ws.Range('D31').Value = '5.97'
$ws.Range('E31').Value = '  +2.39%  '
$ws.Range('D32').Value = '153.35'
$ws.Range('E32').Value = '  +2.17%  '
$ws.Range('D33').Value = '18.93'
$ws.Range('E33').Value = '  -0.34%  '
$ws.Range('D34').Value = '3.90'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('D35').Value = '0.848'
$ws.Range('E35').Value = '  +2.56%  '
$ws.Range('E36').Value = '  -4.57%  '
$ws.Range('D37').Value = '0.819'
$ws.Range('E37').Value = '  -2.44%  '
$ws.Range('D38').Value = '1.42'
$ws.Range('E38').Value = '  -2.97%  '
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('D40').Value = '278.30'
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').Value = '0.589'
$ws.Range('E42').Value = '  -2.71%  '
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').Value = '0.0531'
$ws.Range('E44').Value = '  -2.18%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = '0.0941'
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('D46').Value = '18.47'
$ws.Range('E46').Value = '  -4.56%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '1.902.94'
$ws.Range('E48').Value = '  -3.40%  '
$ws.Range('D49').Value = '17.85'
$ws.Range('E49').Value = '  -2.58%  '
$ws.Range('D50').Value = '4.39'
$ws.Range('E50').Value = '  -3.45%  '
$ws.Range('D51').Value = '108.77'
$ws.Range('E51').Value = '  -1.89%  '
